# Auto-generated COM-interop edit script
# Updates "Recommandations" sheet data (rows 2-48) and "Top_YTD" sheet data (rows 2-10),
# matching the refreshed BRVM recommendation figures from the automated data pull.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Recommandations: the source table now has 47 data rows (was 49); drop the two trailing rows ---
$ws1.Rows.Item(49).Delete()
$ws1.Rows.Item(49).Delete()

# --- Recommandations: rewrite every data row (2-48) with the refreshed values ---
$ws1.Cells.Item(2,1).Value2 = 'BRVM - SERVICES PUBLICS'
$ws1.Cells.Item(2,2).Value2 = 0.0
$ws1.Cells.Item(2,3).Value2 = 8.0
$ws1.Cells.Item(2,4).Value2 = 3232.84
$ws1.Cells.Item(2,5).Value2 = 102.31
$ws1.Cells.Item(2,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(2,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(3,1).Value2 = 'SAFCA CI'
$ws1.Cells.Item(3,2).Value2 = 0.0
$ws1.Cells.Item(3,3).Value2 = 3.0
$ws1.Cells.Item(3,4).Value2 = 2765.0
$ws1.Cells.Item(3,5).Value2 = 985.0
$ws1.Cells.Item(3,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(3,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(4,1).Value2 = 'CFAO MOTORS CI'
$ws1.Cells.Item(4,2).Value2 = 0.0
$ws1.Cells.Item(4,3).Value2 = 4.0
$ws1.Cells.Item(4,4).Value2 = 2600.0
$ws1.Cells.Item(4,5).Value2 = 645.0
$ws1.Cells.Item(4,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(4,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(5,1).Value2 = 'BRVM - AUTRES SECTEURS'
$ws1.Cells.Item(5,2).Value2 = 0.0
$ws1.Cells.Item(5,3).Value2 = 4.0
$ws1.Cells.Item(5,4).Value2 = 2566.13
$ws1.Cells.Item(5,5).Value2 = 641.57
$ws1.Cells.Item(5,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(5,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(6,1).Value2 = 'NEI-CEDA CI'
$ws1.Cells.Item(6,2).Value2 = 0.0
$ws1.Cells.Item(6,3).Value2 = 4.0
$ws1.Cells.Item(6,4).Value2 = 2375.0
$ws1.Cells.Item(6,5).Value2 = 595.0
$ws1.Cells.Item(6,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(6,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(7,1).Value2 = 'SETAO CI'
$ws1.Cells.Item(7,2).Value2 = 0.0
$ws1.Cells.Item(7,3).Value2 = 4.0
$ws1.Cells.Item(7,4).Value2 = 2340.0
$ws1.Cells.Item(7,5).Value2 = 580.0
$ws1.Cells.Item(7,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(7,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(8,1).Value2 = 'UNIWAX CI'
$ws1.Cells.Item(8,2).Value2 = 0.0
$ws1.Cells.Item(8,3).Value2 = 4.0
$ws1.Cells.Item(8,4).Value2 = 2300.0
$ws1.Cells.Item(8,5).Value2 = 580.0
$ws1.Cells.Item(8,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(8,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(9,1).Value2 = 'AIR LIQUIDE CI'
$ws1.Cells.Item(9,2).Value2 = 0.0
$ws1.Cells.Item(9,3).Value2 = 4.0
$ws1.Cells.Item(9,4).Value2 = 2115.0
$ws1.Cells.Item(9,5).Value2 = 545.0
$ws1.Cells.Item(9,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(9,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(10,1).Value2 = 'BRVM - DISTRIBUTION'
$ws1.Cells.Item(10,2).Value2 = 0.0
$ws1.Cells.Item(10,3).Value2 = 4.0
$ws1.Cells.Item(10,4).Value2 = 1452.1
$ws1.Cells.Item(10,5).Value2 = 364.57
$ws1.Cells.Item(10,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(10,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(11,1).Value2 = 'BRVM - TRANSPORT'
$ws1.Cells.Item(11,2).Value2 = 0.0
$ws1.Cells.Item(11,3).Value2 = 4.0
$ws1.Cells.Item(11,4).Value2 = 1402.54
$ws1.Cells.Item(11,5).Value2 = 352.47
$ws1.Cells.Item(11,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(11,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(12,1).Value2 = 'BRVM - AGRICULTURE'
$ws1.Cells.Item(12,2).Value2 = 0.0
$ws1.Cells.Item(12,3).Value2 = 4.0
$ws1.Cells.Item(12,4).Value2 = 1271.51
$ws1.Cells.Item(12,5).Value2 = 321.44
$ws1.Cells.Item(12,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(12,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(13,1).Value2 = 'BRVM - INDUSTRIE'
$ws1.Cells.Item(13,2).Value2 = 0.0
$ws1.Cells.Item(13,3).Value2 = 4.0
$ws1.Cells.Item(13,4).Value2 = 1047.56
$ws1.Cells.Item(13,5).Value2 = 262.61
$ws1.Cells.Item(13,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(13,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(14,1).Value2 = 'BRVM - CONSOMMATION DE BASE'
$ws1.Cells.Item(14,2).Value2 = 0.0
$ws1.Cells.Item(14,3).Value2 = 4.0
$ws1.Cells.Item(14,4).Value2 = 863.29
$ws1.Cells.Item(14,5).Value2 = 216.86
$ws1.Cells.Item(14,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(14,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(15,1).Value2 = 'BRVM-PRINCIPAL'
$ws1.Cells.Item(15,2).Value2 = 0.0
$ws1.Cells.Item(15,3).Value2 = 4.0
$ws1.Cells.Item(15,4).Value2 = 752.94
$ws1.Cells.Item(15,5).Value2 = 188.77
$ws1.Cells.Item(15,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(15,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(16,1).Value2 = 'BRVM - INDUSTRIELS'
$ws1.Cells.Item(16,2).Value2 = 0.0
$ws1.Cells.Item(16,3).Value2 = 4.0
$ws1.Cells.Item(16,4).Value2 = 552.56
$ws1.Cells.Item(16,5).Value2 = 138.37
$ws1.Cells.Item(16,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(16,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(17,1).Value2 = 'BRVM-PRESTIGE'
$ws1.Cells.Item(17,2).Value2 = 0.0
$ws1.Cells.Item(17,3).Value2 = 4.0
$ws1.Cells.Item(17,4).Value2 = 515.41
$ws1.Cells.Item(17,5).Value2 = 128.94
$ws1.Cells.Item(17,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(17,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(18,1).Value2 = 'BRVM - FINANCES'
$ws1.Cells.Item(18,2).Value2 = 0.0
$ws1.Cells.Item(18,3).Value2 = 4.0
$ws1.Cells.Item(18,4).Value2 = 485.23
$ws1.Cells.Item(18,5).Value2 = 121.47
$ws1.Cells.Item(18,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(18,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(19,1).Value2 = 'BRVM - SERVICES FINANCIERS'
$ws1.Cells.Item(19,2).Value2 = 0.0
$ws1.Cells.Item(19,3).Value2 = 4.0
$ws1.Cells.Item(19,4).Value2 = 476.87
$ws1.Cells.Item(19,5).Value2 = 119.38
$ws1.Cells.Item(19,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(19,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(20,1).Value2 = 'BRVM - ENERGIE'
$ws1.Cells.Item(20,2).Value2 = 0.0
$ws1.Cells.Item(20,3).Value2 = 4.0
$ws1.Cells.Item(20,4).Value2 = 430.84
$ws1.Cells.Item(20,5).Value2 = 108.38
$ws1.Cells.Item(20,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(20,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(21,1).Value2 = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$ws1.Cells.Item(21,2).Value2 = 0.0
$ws1.Cells.Item(21,3).Value2 = 4.0
$ws1.Cells.Item(21,4).Value2 = 418.08
$ws1.Cells.Item(21,5).Value2 = 104.7
$ws1.Cells.Item(21,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(21,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(22,1).Value2 = 'BRVM - TELECOMMUNICATIONS'
$ws1.Cells.Item(22,2).Value2 = 0.0
$ws1.Cells.Item(22,3).Value2 = 4.0
$ws1.Cells.Item(22,4).Value2 = 368.31
$ws1.Cells.Item(22,5).Value2 = 91.84
$ws1.Cells.Item(22,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(22,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(23,1).Value2 = 'SAFCA CI (SAFC)'
$ws1.Cells.Item(23,2).Value2 = 4.0
$ws1.Cells.Item(23,3).Value2 = 0.0
$ws1.Cells.Item(23,4).Value2 = 28.66
$ws1.Cells.Item(23,5).Value2 = 7.11
$ws1.Cells.Item(23,6).Value2 = '🟢 Achat'
$ws1.Cells.Item(23,7).Value2 = '✅ Renforcer'
$ws1.Cells.Item(24,1).Value2 = 'SICABLE CI (CABC)'
$ws1.Cells.Item(24,2).Value2 = 2.0
$ws1.Cells.Item(24,3).Value2 = 0.0
$ws1.Cells.Item(24,4).Value2 = 8.17
$ws1.Cells.Item(24,5).Value2 = 4.0
$ws1.Cells.Item(24,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(24,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(25,1).Value2 = 'UNILEVER CI (UNLC)'
$ws1.Cells.Item(25,2).Value2 = 1.0
$ws1.Cells.Item(25,3).Value2 = 0.0
$ws1.Cells.Item(25,4).Value2 = 7.49
$ws1.Cells.Item(25,5).Value2 = 7.49
$ws1.Cells.Item(25,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(25,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(26,1).Value2 = 'BANK OF AFRICA ML (BOAM)'
$ws1.Cells.Item(26,2).Value2 = 1.0
$ws1.Cells.Item(26,3).Value2 = 0.0
$ws1.Cells.Item(26,4).Value2 = 7.37
$ws1.Cells.Item(26,5).Value2 = 7.37
$ws1.Cells.Item(26,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(26,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(27,1).Value2 = 'SOCIETE GENERALE COTE D''IVOIRE (SGBC)'
$ws1.Cells.Item(27,2).Value2 = 1.0
$ws1.Cells.Item(27,3).Value2 = 0.0
$ws1.Cells.Item(27,4).Value2 = 7.04
$ws1.Cells.Item(27,5).Value2 = 7.04
$ws1.Cells.Item(27,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(27,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(28,1).Value2 = 'SAPH CI (SPHC)'
$ws1.Cells.Item(28,2).Value2 = 1.0
$ws1.Cells.Item(28,3).Value2 = 0.0
$ws1.Cells.Item(28,4).Value2 = 4.69
$ws1.Cells.Item(28,5).Value2 = 4.69
$ws1.Cells.Item(28,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(28,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(29,1).Value2 = 'VIVO ENERGY CI (SHEC)'
$ws1.Cells.Item(29,2).Value2 = 1.0
$ws1.Cells.Item(29,3).Value2 = 0.0
$ws1.Cells.Item(29,4).Value2 = 3.35
$ws1.Cells.Item(29,5).Value2 = 3.35
$ws1.Cells.Item(29,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(29,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(30,1).Value2 = 'BERNABE CI (BNBC)'
$ws1.Cells.Item(30,2).Value2 = 1.0
$ws1.Cells.Item(30,3).Value2 = 1.0
$ws1.Cells.Item(30,4).Value2 = 3.35
$ws1.Cells.Item(30,5).Value2 = 7.2
$ws1.Cells.Item(30,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(30,7).Value2 = '👀 À surveiller'
$ws1.Cells.Item(31,1).Value2 = 'AIR LIQUIDE CI (SIVC)'
$ws1.Cells.Item(31,2).Value2 = 1.0
$ws1.Cells.Item(31,3).Value2 = 1.0
$ws1.Cells.Item(31,4).Value2 = 2.94
$ws1.Cells.Item(31,5).Value2 = -1.82
$ws1.Cells.Item(31,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(31,7).Value2 = '👀 À surveiller'
$ws1.Cells.Item(32,1).Value2 = 'SICOR CI (SICC)'
$ws1.Cells.Item(32,2).Value2 = 1.0
$ws1.Cells.Item(32,3).Value2 = 0.0
$ws1.Cells.Item(32,4).Value2 = 2.8
$ws1.Cells.Item(32,5).Value2 = 2.8
$ws1.Cells.Item(32,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(32,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(33,1).Value2 = 'UNIWAX CI (UNXC)'
$ws1.Cells.Item(33,2).Value2 = 1.0
$ws1.Cells.Item(33,3).Value2 = 0.0
$ws1.Cells.Item(33,4).Value2 = 1.75
$ws1.Cells.Item(33,5).Value2 = 1.75
$ws1.Cells.Item(33,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(33,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(34,1).Value2 = 'NEI-CEDA CI (NEIC)'
$ws1.Cells.Item(34,2).Value2 = 1.0
$ws1.Cells.Item(34,3).Value2 = 1.0
$ws1.Cells.Item(34,4).Value2 = 0.85
$ws1.Cells.Item(34,5).Value2 = 1.69
$ws1.Cells.Item(34,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(34,7).Value2 = '👀 À surveiller'
$ws1.Cells.Item(35,1).Value2 = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws1.Cells.Item(35,2).Value2 = 1.0
$ws1.Cells.Item(35,3).Value2 = 1.0
$ws1.Cells.Item(35,4).Value2 = 0.42
$ws1.Cells.Item(35,5).Value2 = 6.67
$ws1.Cells.Item(35,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(35,7).Value2 = '👀 À surveiller'
$ws1.Cells.Item(36,1).Value2 = 'BANK OF AFRICA NG (BOAN)'
$ws1.Cells.Item(36,2).Value2 = 1.0
$ws1.Cells.Item(36,3).Value2 = 1.0
$ws1.Cells.Item(36,4).Value2 = 0.05
$ws1.Cells.Item(36,5).Value2 = 2.24
$ws1.Cells.Item(36,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(36,7).Value2 = '👀 À surveiller'
$ws1.Cells.Item(37,1).Value2 = 'SETAO CI (STAC)'
$ws1.Cells.Item(37,2).Value2 = 1.0
$ws1.Cells.Item(37,3).Value2 = 1.0
$ws1.Cells.Item(37,4).Value2 = 0.03
$ws1.Cells.Item(37,5).Value2 = -3.45
$ws1.Cells.Item(37,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(37,7).Value2 = '👀 À surveiller'
$ws1.Cells.Item(38,1).Value2 = 'TOTAL'
$ws1.Cells.Item(38,2).Value2 = 0.0
$ws1.Cells.Item(38,3).Value2 = 4.0
$ws1.Cells.Item(38,4).Value2 = 0.0
$ws1.Cells.Item(38,5).Value2 = 0.0
$ws1.Cells.Item(38,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(38,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(39,1).Value2 = 'FILTISAC CI (FTSC)'
$ws1.Cells.Item(39,2).Value2 = 0.0
$ws1.Cells.Item(39,3).Value2 = 1.0
$ws1.Cells.Item(39,4).Value2 = -0.8
$ws1.Cells.Item(39,5).Value2 = -0.8
$ws1.Cells.Item(39,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(39,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(40,1).Value2 = 'SOLIBRA CI (SLBC)'
$ws1.Cells.Item(40,2).Value2 = 0.0
$ws1.Cells.Item(40,3).Value2 = 1.0
$ws1.Cells.Item(40,4).Value2 = -0.81
$ws1.Cells.Item(40,5).Value2 = -0.81
$ws1.Cells.Item(40,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(40,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(41,1).Value2 = 'ORANGE COTE D''IVOIRE (ORAC)'
$ws1.Cells.Item(41,2).Value2 = 0.0
$ws1.Cells.Item(41,3).Value2 = 1.0
$ws1.Cells.Item(41,4).Value2 = -1.06
$ws1.Cells.Item(41,5).Value2 = -1.06
$ws1.Cells.Item(41,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(41,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(42,1).Value2 = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$ws1.Cells.Item(42,2).Value2 = 0.0
$ws1.Cells.Item(42,3).Value2 = 1.0
$ws1.Cells.Item(42,4).Value2 = -1.32
$ws1.Cells.Item(42,5).Value2 = -1.32
$ws1.Cells.Item(42,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(42,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(43,1).Value2 = 'SONATEL SN (SNTS)'
$ws1.Cells.Item(43,2).Value2 = 0.0
$ws1.Cells.Item(43,3).Value2 = 1.0
$ws1.Cells.Item(43,4).Value2 = -1.4
$ws1.Cells.Item(43,5).Value2 = -1.4
$ws1.Cells.Item(43,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(43,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(44,1).Value2 = 'LOTERIE NATIONALE DU BENIN (LNBB)'
$ws1.Cells.Item(44,2).Value2 = 1.0
$ws1.Cells.Item(44,3).Value2 = 1.0
$ws1.Cells.Item(44,4).Value2 = -1.58
$ws1.Cells.Item(44,5).Value2 = 5.88
$ws1.Cells.Item(44,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(44,7).Value2 = '👀 À surveiller'
$ws1.Cells.Item(45,1).Value2 = 'BICI CI (BICC)'
$ws1.Cells.Item(45,2).Value2 = 0.0
$ws1.Cells.Item(45,3).Value2 = 1.0
$ws1.Cells.Item(45,4).Value2 = -2.85
$ws1.Cells.Item(45,5).Value2 = -2.85
$ws1.Cells.Item(45,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(45,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(46,1).Value2 = 'ONATEL BF (ONTBF)'
$ws1.Cells.Item(46,2).Value2 = 0.0
$ws1.Cells.Item(46,3).Value2 = 2.0
$ws1.Cells.Item(46,4).Value2 = -4.21
$ws1.Cells.Item(46,5).Value2 = -2.13
$ws1.Cells.Item(46,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(46,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(47,1).Value2 = 'ORAGROUP TOGO (ORGT)'
$ws1.Cells.Item(47,2).Value2 = 0.0
$ws1.Cells.Item(47,3).Value2 = 2.0
$ws1.Cells.Item(47,4).Value2 = -5.21
$ws1.Cells.Item(47,5).Value2 = -1.57
$ws1.Cells.Item(47,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(47,7).Value2 = '➖ Neutre'
$ws1.Cells.Item(48,1).Value2 = 'CIE CI (CIEC)'
$ws1.Cells.Item(48,2).Value2 = 0.0
$ws1.Cells.Item(48,3).Value2 = 2.0
$ws1.Cells.Item(48,4).Value2 = -6.62
$ws1.Cells.Item(48,5).Value2 = -2.71
$ws1.Cells.Item(48,6).Value2 = '🟡 Observer'
$ws1.Cells.Item(48,7).Value2 = '➖ Neutre'

# --- Top_YTD: rewrite every data row (2-10) with the refreshed values ---
$ws2.Cells.Item(2,1).Value2 = 'BRVM - SERVICES PUBLICS'
$ws2.Cells.Item(2,2).Value2 = 7059160.66
$ws2.Cells.Item(3,1).Value2 = 'CFAO MOTORS CI'
$ws2.Cells.Item(3,2).Value2 = 316292.19
$ws2.Cells.Item(4,1).Value2 = 'BRVM - AUTRES SECTEURS'
$ws2.Cells.Item(4,2).Value2 = 301851.9
$ws2.Cells.Item(5,1).Value2 = 'NEI-CEDA CI'
$ws2.Cells.Item(5,2).Value2 = 231534.64
$ws2.Cells.Item(6,1).Value2 = 'SETAO CI'
$ws2.Cells.Item(6,2).Value2 = 220048.64
$ws2.Cells.Item(7,1).Value2 = 'UNIWAX CI'
$ws2.Cells.Item(7,2).Value2 = 207482.75
$ws2.Cells.Item(8,1).Value2 = 'AIR LIQUIDE CI'
$ws2.Cells.Item(8,2).Value2 = 156110.94
$ws2.Cells.Item(9,1).Value2 = 'SAFCA CI'
$ws2.Cells.Item(9,2).Value2 = 106143.2
$ws2.Cells.Item(10,1).Value2 = 'BRVM - DISTRIBUTION'
$ws2.Cells.Item(10,2).Value2 = 45863.6
$ws2.Cells.Item(11,1).Value2 = 'BRVM - TRANSPORT'
$ws2.Cells.Item(11,2).Value2 = 41134.07
